# Scheduled-runner data refresh: update computed columns H:N (market-price-derived
# leve profit figures) for a handful of rows across the per-job sheets. No formulas
# are involved -- these are plain cached values that get overwritten in place, and
# a couple of rows gain/lose an H..N cell outright (the source value was blank/absent).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused
$ws.Range("H33").Value = 8404.379000000001
$ws.Range("I33").Value = 11012.286
$ws.Range("J33").Value = 1558.625
$ws.Range("K33").Value = 11012.286
$ws.Range("L33").Value = 1558.625
$ws.Range("M33").Value = -10783.286
$ws.Range("N33").Value = -2016.625
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 2000.25
$ws.Range("I62").Value = 2000.25
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2000.25
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -1376.25
# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 2000.25
$ws.Range("I65").Value = 2000.25
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 10001.25
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -6881.25
# Row 98: The Dotted Line
$ws.Range("H98").Value = 741.1818
$ws.Range("I98").Value = 615.3
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 615.3
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 882.7
$ws.Range("N98").Value = -4996
# Row 122: Wishful Inking
$ws.Range("H122").Value = 741.1818
$ws.Range("I122").Value = 615.3
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 1845.9
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = 604.1000000000001
$ws.Range("N122").Value = -10900
# Row 138: All-night Crafting
$ws.Range("H138").Value = 16396552
$ws.Range("I138").Value = 1333.3334
$ws.Range("J138").Value = 27031288
$ws.Range("K138").Value = 4000.0002
$ws.Range("L138").Value = 81093864
$ws.Range("M138").Value = 1139.9998
$ws.Range("N138").Value = -81104144

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 3766.2922
$ws.Range("I32").Value = 3870.459
$ws.Range("J32").Value = 2177.75
$ws.Range("K32").Value = 3870.459
$ws.Range("L32").Value = 2177.75
$ws.Range("M32").Value = -3583.459
$ws.Range("N32").Value = -2751.75
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 161535.08
$ws.Range("I74").Value = 188125.25
$ws.Range("J74").Value = 1994
$ws.Range("K74").Value = 188125.25
$ws.Range("L74").Value = 1994
$ws.Range("M74").Value = -187251.25
$ws.Range("N74").Value = -3742
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 161535.08
$ws.Range("I77").Value = 188125.25
$ws.Range("J77").Value = 1994
$ws.Range("K77").Value = 940626.25
$ws.Range("L77").Value = 9970
$ws.Range("M77").Value = -936258.25
$ws.Range("N77").Value = -18706

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 1838.5238
$ws.Range("I134").Value = 1329.9412
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 3989.8236
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -1454.8236
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 5753.909
$ws.Range("I31").Value = 4385.5625
$ws.Range("J31").Value = 7041.7646
$ws.Range("K31").Value = 4385.5625
$ws.Range("L31").Value = 7041.7646
$ws.Range("M31").Value = -4090.5625
$ws.Range("N31").Value = -7631.7646
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 5753.909
$ws.Range("I34").Value = 4385.5625
$ws.Range("J34").Value = 7041.7646
$ws.Range("K34").Value = 4385.5625
$ws.Range("L34").Value = 7041.7646
$ws.Range("M34").Value = -4183.5625
$ws.Range("N34").Value = -7445.7646
# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 3350.875
$ws.Range("I122").Value = 2252
$ws.Range("J122").Value = 4449.75
$ws.Range("K122").Value = 6756
$ws.Range("L122").Value = 13349.25
$ws.Range("M122").Value = -4306
$ws.Range("N122").Value = -18249.25
# Row 131: An Integral Reward
$ws.Range("H131").Value = 29833.875
$ws.Range("I131").Value = 12000
$ws.Range("J131").Value = 32381.572
$ws.Range("K131").Value = 12000
$ws.Range("L131").Value = 32381.572
$ws.Range("M131").Value = -6960
$ws.Range("N131").Value = -42461.572
# Row 135: The Wing's Wings
$ws.Range("H135").Value = 120471.875
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 120471.875
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 120471.875
$ws.Range("N135").Value = -130611.875
# Row 141: No Greater Treasure
$ws.Range("H141").Value = 97360.086
$ws.Range("I141").Value = 31250
$ws.Range("J141").Value = 130415.125
$ws.Range("K141").Value = 31250
$ws.Range("L141").Value = 130415.125
$ws.Range("M141").Value = -26070
$ws.Range("N141").Value = -140775.125

$ws = $wb.Worksheets.Item("CUL")
# Row 14: Keep Your Powder Dry
$ws.Range("H14").Value = 207.38461
$ws.Range("I14").Value = 207.38461
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 622.15383
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -449.15383
# Row 38: Pretty as a Picture
$ws.Range("H38").Value = 33.25
$ws.Range("I38").Value = 15.25
$ws.Range("J38").Value = 51.25
$ws.Range("K38").Value = 45.75
$ws.Range("L38").Value = 153.75
$ws.Range("M38").Value = 301.25
$ws.Range("N38").Value = -847.75
# Row 48: Rise and Dine
$ws.Range("H48").Value = 2067.2727
$ws.Range("I48").Value = 677.8570999999999
$ws.Range("J48").Value = 4498.75
$ws.Range("K48").Value = 2033.5713
$ws.Range("L48").Value = 13496.25
$ws.Range("M48").Value = -1783.5713
$ws.Range("N48").Value = -13996.25
# Row 68: Such a Butter Face
$ws.Range("H68").Value = 1499.6666
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1499.6666
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 4498.9998
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -6120.9998
# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 1499.6666
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1499.6666
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 13496.9994
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -21608.9994
# Row 109: Cure for What Ails
$ws.Range("H109").Value = 6949.75
$ws.Range("I109").Value = 8719.6
$ws.Range("J109").Value = 4000
$ws.Range("K109").Value = 26158.8
$ws.Range("L109").Value = 12000
$ws.Range("M109").Value = -25118.8
$ws.Range("N109").Value = -14080
# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 3432
$ws.Range("I113").Value = 995
$ws.Range("J113").Value = 3736.625
$ws.Range("K113").Value = 2985
$ws.Range("L113").Value = 11209.875
$ws.Range("M113").Value = -815
$ws.Range("N113").Value = -15549.875
# Row 136: Simple Is Hardest
$ws.Range("H136").Value = 4252.5
$ws.Range("I136").Value = 3505
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 10515
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -5415
$ws.Range("N136").Value = -25200

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 76656.67999999999
$ws.Range("I70").Value = 99344.48
$ws.Range("J70").Value = 8593.286
$ws.Range("K70").Value = 99344.48
$ws.Range("L70").Value = 8593.286
$ws.Range("M70").Value = -99074.48
$ws.Range("N70").Value = -9133.286
# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 76656.67999999999
$ws.Range("I73").Value = 99344.48
$ws.Range("J73").Value = 8593.286
$ws.Range("K73").Value = 99344.48
$ws.Range("L73").Value = 8593.286
$ws.Range("M73").Value = -98408.48
$ws.Range("N73").Value = -10465.286
# Row 132: On Board for Lar
$ws.Range("H132").Value = 1540.9
$ws.Range("I132").Value = 1358.5714
$ws.Range("J132").Value = 1966.3334
$ws.Range("K132").Value = 4075.7142
$ws.Range("L132").Value = 5899.0002
$ws.Range("M132").Value = -1545.7142
$ws.Range("N132").Value = -10959.0002

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 8477
$ws.Range("I7").Value = 8451.5
$ws.Range("J7").Value = 8502.5
$ws.Range("K7").Value = 8451.5
$ws.Range("L7").Value = 8502.5
$ws.Range("M7").Value = -8339.5
$ws.Range("N7").Value = -8726.5
# Row 16: Saddle Sore
$ws.Range("H16").Value = 5001.1113
$ws.Range("I16").Value = 5000
$ws.Range("J16").Value = 5001.6665
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 5001.6665
$ws.Range("M16").Value = -4830
$ws.Range("N16").Value = -5341.6665
# Row 122: Hell on Leather
$ws.Range("H122").Value = 3754.9333
$ws.Range("I122").Value = 3477.4167
$ws.Range("J122").Value = 4865
$ws.Range("K122").Value = 10432.2501
$ws.Range("L122").Value = 14595
$ws.Range("M122").Value = -7982.250100000001
$ws.Range("N122").Value = -19495
# Row 126: Battered Books
$ws.Range("H126").Value = 8477
$ws.Range("I126").Value = 8451.5
$ws.Range("J126").Value = 8502.5
$ws.Range("K126").Value = 25354.5
$ws.Range("L126").Value = 25507.5
$ws.Range("M126").Value = -22884.5
$ws.Range("N126").Value = -30447.5
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 8751.3125
$ws.Range("I132").Value = 2578.5386
$ws.Range("J132").Value = 35500
$ws.Range("K132").Value = 7735.6158
$ws.Range("L132").Value = 106500
$ws.Range("M132").Value = -5205.6158
$ws.Range("N132").Value = -111560
# Row 134: Freezing Fingers
$ws.Range("H134").Value = 135214
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 135214
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 135214
$ws.Range("N134").Value = -145354
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 2155.2646
$ws.Range("I136").Value = 1771.4138
$ws.Range("J136").Value = 4381.6
$ws.Range("K136").Value = 5314.2414
$ws.Range("L136").Value = 13144.8
$ws.Range("M136").Value = -2764.2414
$ws.Range("N136").Value = -18244.8

$ws = $wb.Worksheets.Item("WVR")
# Row 17: Making Gloves Out of Nothing at All
$ws.Range("H17").Value = 336
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 4
$ws.Range("L17").Value = 1000
$ws.Range("M17").Value = 168
$ws.Range("N17").Value = -1344
# Row 27: Hitting Below the Belt
$ws.Range("H27").Value = 79992
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 79992
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 79992
$ws.Range("N27").Value = -80130
# Row 46: Crunching the Numbers
$ws.Range("H46").Value = 154999.67
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 154999.67
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 154999.67
$ws.Range("N46").Value = -155461.67
# Row 115: Gloves Come in Handy
$ws.Range("H115").Value = 63330.668
$ws.Range("I115").Value = 29998
$ws.Range("J115").Value = 79997
$ws.Range("K115").Value = 29998
$ws.Range("L115").Value = 79997
$ws.Range("M115").Value = -28431
$ws.Range("N115").Value = -83131
# Row 119: A Job Well Done
$ws.Range("H119").Value = 84995
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 84995
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 84995
$ws.Range("N119").Value = -94671
# Row 134: Cloth for Canvas
$ws.Range("H134").Value = 154999.67
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 154999.67
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 464999.01
$ws.Range("N134").Value = -470069.01
